$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header columns: Original Facility* / Destination Facility* ---
$ws.Range("H1").Value = "Original Facility*"
$ws.Range("I1").Value = "Destination Facility*"

# --- Row 2 updates ---
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = 44639
$ws.Range("H2").Value = "test11"
$ws.Range("I2").Value = "test12"

# --- Row 3: clear everything except the (still date-styled) B3/C3 blanks ---
$ws.Range("A3").ClearContents()
$ws.Range("B3").ClearContents()
$ws.Range("D3:G3").ClearContents()

# --- Column widths for the two new columns (best-effort match to target bestFit widths) ---
$ws.Columns.Item(8).ColumnWidth = 14.9
$ws.Columns.Item(9).ColumnWidth = 18.3

# --- Selection moves to B2 ---
$ws.Range("B2").Select()
